$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.363.30"
$ws.Range("E2").Value = "  -4.44%  "

$ws.Range("D3").Value = "1.569.63"
$ws.Range("E3").Value = "  -4.10%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3688"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.55%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.16"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.24%  "

$ws.Range("E9").Value = "  -4.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.164"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07599"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.046"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.892"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.07%  "

$ws.Range("D16").Value = "1.569.14"
$ws.Range("E16").Value = "  -4.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001132"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.94%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06753"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.78%  "

$ws.Range("E20").Value = "  -0.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.241"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.01%  "

$ws.Range("E22").Value = "  -4.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5314"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.12%  "

$ws.Range("D25").Value = "22.379.43"
$ws.Range("E25").Value = "  -4.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.383"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.974"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "145.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.957"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.74%  "

$ws.Range("D32").Value = "1.747.17"
$ws.Range("E32").Value = "  -3.99%  "

$ws.Range("E33").Value = "  +6.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.250"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.996"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.34%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08445"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02528"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2328"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.545"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06501"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.241"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6361"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9996"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5973"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.752"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.89%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.125"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.251"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "123.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.56%  "
